$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.59"
$ws.Range("D4").Value = "'5.349"
$ws.Range("D5").Value = "'0.05635"
$ws.Range("D6").Value = "'3.428"
$ws.Range("D7").Value = "'6.368"
$ws.Range("D8").Value = "'0.8187"
$ws.Range("D9").Value = "'0.9331"
$ws.Range("D10").Value = "'0.1441"
$ws.Range("D11").Value = "'0.07487"
$ws.Range("D12").Value = "'0.03255"
$ws.Range("D13").Value = "'0.03094"
$ws.Range("D14").Value = "'0.09303"
$ws.Range("D15").Value = "'3.564"
$ws.Range("D16").Value = "'0.001636"
$ws.Range("D18").Value = "'0.0005780"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006325"
$ws.Range("D20").Value = "'0.005060"
$ws.Range("D21").Value = "'0.001033"
$ws.Range("D22").Value = "'0.0001499"
$ws.Range("D23").Value = "'3.755"
$ws.Range("D25").Value = "'0.3307"
$ws.Range("D26").Value = "'0.1328"
$ws.Range("D28").Value = "'0.0003000"
$ws.Range("D40").Value = "'0.03952"
$ws.Range("D41").Value = "'0.002915"
$ws.Range("E41").Value = "40KickTokenKICKWorstin24h"
$ws.Range("D43").Value = "'0.003019"
$ws.Range("D44").Value = "'0.008582"
$ws.Range("D45").Value = "'0.00005580"
$ws.Range("D47").Value = "'0.0005500"
$ws.Range("D48").Value = "'0.7800"
$ws.Range("D49").Value = "'0.1769"
